# Auto-generated edit script: updates crypto Price (D) and Volume(1h) (E) columns
# to reflect the refreshed values from the GitHub Actions crypto-list update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.731.25"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "'3.300.16"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'187.34"
$ws.Range("E5").Value = "  +5.66%  "
$ws.Range("D6").Value = "'552.64"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.578"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "'3.290.52"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'46.59"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "'8.59"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "'3.826.76"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "'594.90"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "'65.709.68"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "'17.81"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'3.297.88"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'10.96"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  +7.32%  "
$ws.Range("D24").Value = "'5.06"
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").Value = "'100.46"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").Value = "'9.44"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").Value = "'8.64"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "'30.12"
$ws.Range("D32").Value = "'6.67"
$ws.Range("E32").Value = "  +8.70%  "
$ws.Range("D33").Value = "'3.82"
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("D34").Value = "'568.37"
$ws.Range("E34").Value = "  +8.47%  "
$ws.Range("D35").Value = "'10.96"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'56.87"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").Value = "'3.689.81"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'3.44"
$ws.Range("E40").Value = "  +10.06%  "
$ws.Range("D41").Value = "'33.51"
$ws.Range("E41").Value = "  +7.07%  "
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").Value = "'0.0₃0688"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").Value = "'3.40"
$ws.Range("E46").Value = "  +6.11%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  -0.03%  "
